$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns A (dates) and F (mixed numeric/text) to Text format
# so values are stored as literal strings, matching the source data
# (all cells in this sheet are plain text, not numbers/dates).
$ws.Range("A158:A183").NumberFormat = "@"
$ws.Range("F158:F183").NumberFormat = "@"

$ws.Range("A158").Value = "2024-07-25"
$ws.Range("B158").Value = "13:15:09"
$ws.Range("C158").Value = "Disney's Hollywood Studios"
$ws.Range("D158").Value = "Toy Story Land"
$ws.Range("E158").Value = "Alien Swirling Saucers"
$ws.Range("F158").Value = "40"

$ws.Range("A159").Value = "2024-07-25"
$ws.Range("B159").Value = "13:15:09"
$ws.Range("C159").Value = "Disney's Hollywood Studios"
$ws.Range("D159").Value = "Commissary Lane"
$ws.Range("E159").Value = "Meet Disney Stars at Red Carpet Dreams"
$ws.Range("F159").Value = "40"

$ws.Range("A160").Value = "2024-07-25"
$ws.Range("B160").Value = "13:15:09"
$ws.Range("C160").Value = "Disney's Hollywood Studios"
$ws.Range("D160").Value = "Echo Lake"
$ws.Range("E160").Value = "Meet Olaf at Celebrity Spotlight"
$ws.Range("F160").Value = "25"

$ws.Range("A161").Value = "2024-07-25"
$ws.Range("B161").Value = "13:15:09"
$ws.Range("C161").Value = "Disney's Hollywood Studios"
$ws.Range("D161").Value = "Hollywood Boulevard"
$ws.Range("E161").Value = "Mickey & Minnie's Runaway Railway"
$ws.Range("F161").Value = "65"

$ws.Range("A162").Value = "2024-07-25"
$ws.Range("B162").Value = "13:15:09"
$ws.Range("C162").Value = "Disney's Hollywood Studios"
$ws.Range("D162").Value = "Star Wars: Galaxy's Edge"
$ws.Range("E162").Value = "Millennium Falcon: Smugglers Run"
$ws.Range("F162").Value = "85"

$ws.Range("A163").Value = "2024-07-25"
$ws.Range("B163").Value = "13:15:09"
$ws.Range("C163").Value = "Disney's Hollywood Studios"
$ws.Range("D163").Value = "Grand Avenue"
$ws.Range("E163").Value = "Muppet*Vision 3D"
$ws.Range("F163").Value = "10"

$ws.Range("A164").Value = "2024-07-25"
$ws.Range("B164").Value = "13:15:09"
$ws.Range("C164").Value = "Disney's Hollywood Studios"
$ws.Range("D164").Value = "Sunset Boulevard"
$ws.Range("E164").Value = "Rock 'n' Roller Coaster Starring Aerosmith"
$ws.Range("F164").Value = "95"

$ws.Range("A165").Value = "2024-07-25"
$ws.Range("B165").Value = "13:15:09"
$ws.Range("C165").Value = "Disney's Hollywood Studios"
$ws.Range("D165").Value = "Toy Story Land"
$ws.Range("E165").Value = "Slinky Dog Dash"
$ws.Range("F165").Value = "85"

$ws.Range("A166").Value = "2024-07-25"
$ws.Range("B166").Value = "13:15:09"
$ws.Range("C166").Value = "Disney's Hollywood Studios"
$ws.Range("D166").Value = "Echo Lake"
$ws.Range("E166").Value = "Star Tours – The Adventures Continue"
$ws.Range("F166").Value = "30"

$ws.Range("A167").Value = "2024-07-25"
$ws.Range("B167").Value = "13:15:09"
$ws.Range("C167").Value = "Disney's Hollywood Studios"
$ws.Range("D167").Value = "Animation Courtyard"
$ws.Range("E167").Value = "Star Wars Launch Bay: Meet Chewbacca"
$ws.Range("F167").Value = "40"

$ws.Range("A168").Value = "2024-07-25"
$ws.Range("B168").Value = "13:15:09"
$ws.Range("C168").Value = "Disney's Hollywood Studios"
$ws.Range("D168").Value = "Star Wars: Galaxy's Edge"
$ws.Range("E168").Value = "Star Wars: Rise of the Resistance"
$ws.Range("F168").Value = "Atração indisponível agora"

$ws.Range("A169").Value = "2024-07-25"
$ws.Range("B169").Value = "13:15:09"
$ws.Range("C169").Value = "Disney's Hollywood Studios"
$ws.Range("D169").Value = "Toy Story Land"
$ws.Range("E169").Value = "Toy Story Mania!"
$ws.Range("F169").Value = "60"

$ws.Range("A170").Value = "2024-07-25"
$ws.Range("B170").Value = "13:15:09"
$ws.Range("C170").Value = "Disney's Hollywood Studios"
$ws.Range("D170").Value = "Sunset Boulevard"
$ws.Range("E170").Value = "The Twilight Zone Tower of Terror"
$ws.Range("F170").Value = "120"

$ws.Range("A171").Value = "2024-07-25"
$ws.Range("B171").Value = "13:20:22"
$ws.Range("C171").Value = "Disney's Hollywood Studios"
$ws.Range("D171").Value = "Toy Story Land"
$ws.Range("E171").Value = "Alien Swirling Saucers"
$ws.Range("F171").Value = "40"

$ws.Range("A172").Value = "2024-07-25"
$ws.Range("B172").Value = "13:20:22"
$ws.Range("C172").Value = "Disney's Hollywood Studios"
$ws.Range("D172").Value = "Commissary Lane"
$ws.Range("E172").Value = "Meet Disney Stars at Red Carpet Dreams"
$ws.Range("F172").Value = "40"

$ws.Range("A173").Value = "2024-07-25"
$ws.Range("B173").Value = "13:20:22"
$ws.Range("C173").Value = "Disney's Hollywood Studios"
$ws.Range("D173").Value = "Echo Lake"
$ws.Range("E173").Value = "Meet Olaf at Celebrity Spotlight"
$ws.Range("F173").Value = "25"

$ws.Range("A174").Value = "2024-07-25"
$ws.Range("B174").Value = "13:20:22"
$ws.Range("C174").Value = "Disney's Hollywood Studios"
$ws.Range("D174").Value = "Hollywood Boulevard"
$ws.Range("E174").Value = "Mickey & Minnie's Runaway Railway"
$ws.Range("F174").Value = "65"

$ws.Range("A175").Value = "2024-07-25"
$ws.Range("B175").Value = "13:20:22"
$ws.Range("C175").Value = "Disney's Hollywood Studios"
$ws.Range("D175").Value = "Star Wars: Galaxy's Edge"
$ws.Range("E175").Value = "Millennium Falcon: Smugglers Run"
$ws.Range("F175").Value = "85"

$ws.Range("A176").Value = "2024-07-25"
$ws.Range("B176").Value = "13:20:22"
$ws.Range("C176").Value = "Disney's Hollywood Studios"
$ws.Range("D176").Value = "Grand Avenue"
$ws.Range("E176").Value = "Muppet*Vision 3D"
$ws.Range("F176").Value = "10"

$ws.Range("A177").Value = "2024-07-25"
$ws.Range("B177").Value = "13:20:22"
$ws.Range("C177").Value = "Disney's Hollywood Studios"
$ws.Range("D177").Value = "Sunset Boulevard"
$ws.Range("E177").Value = "Rock 'n' Roller Coaster Starring Aerosmith"
$ws.Range("F177").Value = "30"

$ws.Range("A178").Value = "2024-07-25"
$ws.Range("B178").Value = "13:20:22"
$ws.Range("C178").Value = "Disney's Hollywood Studios"
$ws.Range("D178").Value = "Toy Story Land"
$ws.Range("E178").Value = "Slinky Dog Dash"
$ws.Range("F178").Value = "85"

$ws.Range("A179").Value = "2024-07-25"
$ws.Range("B179").Value = "13:20:22"
$ws.Range("C179").Value = "Disney's Hollywood Studios"
$ws.Range("D179").Value = "Echo Lake"
$ws.Range("E179").Value = "Star Tours – The Adventures Continue"
$ws.Range("F179").Value = "25"

$ws.Range("A180").Value = "2024-07-25"
$ws.Range("B180").Value = "13:20:22"
$ws.Range("C180").Value = "Disney's Hollywood Studios"
$ws.Range("D180").Value = "Animation Courtyard"
$ws.Range("E180").Value = "Star Wars Launch Bay: Meet Chewbacca"
$ws.Range("F180").Value = "40"

$ws.Range("A181").Value = "2024-07-25"
$ws.Range("B181").Value = "13:20:22"
$ws.Range("C181").Value = "Disney's Hollywood Studios"
$ws.Range("D181").Value = "Star Wars: Galaxy's Edge"
$ws.Range("E181").Value = "Star Wars: Rise of the Resistance"
$ws.Range("F181").Value = "75"

$ws.Range("A182").Value = "2024-07-25"
$ws.Range("B182").Value = "13:20:22"
$ws.Range("C182").Value = "Disney's Hollywood Studios"
$ws.Range("D182").Value = "Toy Story Land"
$ws.Range("E182").Value = "Toy Story Mania!"
$ws.Range("F182").Value = "60"

$ws.Range("A183").Value = "2024-07-25"
$ws.Range("B183").Value = "13:20:22"
$ws.Range("C183").Value = "Disney's Hollywood Studios"
$ws.Range("D183").Value = "Sunset Boulevard"
$ws.Range("E183").Value = "The Twilight Zone Tower of Terror"
$ws.Range("F183").Value = "120"
